# Apply the recorded changes from the diff:
# - Row 2: round Q2/R2 to nearest integer; remove Z2/AB2 (Starttid/Sluttid "00:00")
# - Row 3 and Row 4: swap the taxon-specific fields (A,B,E,F,G,H) and set the
#   rounded, row-swapped Ost/Nord (Q/R) coordinates; remove Z3/AB3/Z4/AB4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("Q2").Value2 = 756188
$ws.Range("R2").Value2 = 7291007
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# --- Row 3 (becomes the old Row 4 data) ---
$ws.Range("A3").Value2 = 112181650
$ws.Range("B3").Value2 = 78578
$ws.Range("E3").Value2 = 6458
$ws.Range("F3").Value2 = "Lunglav"
$ws.Range("G3").Value2 = "Lobaria pulmonaria"
$ws.Range("H3").Value2 = "(L.) Hoffm."
$ws.Range("Q3").Value2 = 756202
$ws.Range("R3").Value2 = 7291065
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# --- Row 4 (becomes the old Row 3 data) ---
$ws.Range("A4").Value2 = 112182541
$ws.Range("B4").Value2 = 77267
$ws.Range("E4").Value2 = 6446
$ws.Range("F4").Value2 = "Kolflarnlav"
$ws.Range("G4").Value2 = "Carbonicola anthracophila"
$ws.Range("H4").Value2 = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q4").Value2 = 756204
$ws.Range("R4").Value2 = 7291065
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
